# Daily auto push update: insert a new data row for 2026/02/08 21:00
# just before the existing row 774 (2026/12/29), shifting all subsequent
# rows down by one (old row 815 becomes new row 816).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 774; everything currently at/after
# row 774 shifts down by one row (standard Excel "insert row" behavior).
$ws.Rows.Item(774).Insert()

# Populate the newly inserted row 774 with the new data point.
# The leading apostrophe forces column A to be stored as literal text
# (matching the existing date-as-text convention used throughout column A)
# instead of being auto-converted into a date serial value.
$ws.Range("A774").Value = "'2026/02/08"
$ws.Range("B774").Value = "日"
$ws.Range("C774").Value = 21
$ws.Range("D774").Value = 103
